# Append new work-log entries (rows 28-32) to Sheet1 and update the
# selected/scrolled view to match the new extent of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows: Column B = time span, Column C = work content
$newRows = @(
    @("2019.5.30  18：00-21：00", "对项目整体进行测试和改错"),
    @("2019.5.31  17：00-19：00", "对项目整体进行测试和改错"),
    @("2019.6.5   18：30-20：30", "学习dispatcherservlet配置"),
    @("2019.6.6   18：30-20：00", "学习一个简单的web的helloword例子"),
    @("2019.6.7   16：30-18：30", "配置web层的基本文件")
)

$startRow = 28
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][1]
}

# Update the visible selection / scroll position to reflect the new last cell
$lastRow = $startRow + $newRows.Count - 1
$ws.Range("C" + $lastRow).Select()
$excel.ActiveWindow.ScrollRow = 16
